# Edit the "table" worksheet so that the former root node "Hard Skills"
# becomes "Skills", and "Soft Skills" becomes a child of the new "Skills"
# root (matching the target diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table")

# Rename the "Hard Skills" node (id + label) to "Skills"
$ws.Range("A3").Value = "Skills"
$ws.Range("B3").Value = "Skills"

# "Soft Skills" (row 2) now becomes a child of the new "Skills" root
$ws.Range("C2").Value = "Skills"

# Update the parent references of all former direct children of
# "Hard Skills" (rows 12 through 23) to point at the renamed "Skills" node
for ($r = 12; $r -le 23; $r++) {
    $ws.Cells.Item($r, 3).Value = "Skills"
}

# Restore the active cell selection recorded in the saved workbook
$ws.Range("F13").Select()

$wb.Save()
